$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 127
$ws.Cells.Item(127,1).Value = 126
$ws.Cells.Item(127,2).Value = "Friday, Jan 13"
$ws.Cells.Item(127,3).Value = "3:40 PM"
$ws.Cells.Item(127,4).Value = "FR7894"
$ws.Cells.Item(127,5).Value = "Venice"
$ws.Cells.Item(127,6).Value = "(TSF)"
$ws.Cells.Item(127,7).Value = "Ryanair "
$ws.Cells.Item(127,8).Value = "B738"
$ws.Cells.Item(127,9).Value = "(SP-RSM)"
$ws.Cells.Item(127,10).Value = "3:24 PM"
$ws.Cells.Item(127,11).Borders.LineStyle = 0
$ws.Cells.Item(127,12).Value = "0 hours, -16 minutes"
$ws.Cells.Item(127,13).Borders.LineStyle = 0

# Row 128
$ws.Cells.Item(128,1).Value = 127
$ws.Cells.Item(128,2).Value = "Friday, Jan 13"
$ws.Cells.Item(128,3).Value = "4:00 PM"
$ws.Cells.Item(128,4).Value = "W91902"
$ws.Cells.Item(128,5).Value = "London"
$ws.Cells.Item(128,6).Value = "(LTN)"
$ws.Cells.Item(128,7).Value = "Wizz Air "
$ws.Cells.Item(128,8).Value = "A320"
$ws.Cells.Item(128,9).Value = "(G-WUKF)"
$ws.Cells.Item(128,10).Value = "3:30 PM"
$ws.Cells.Item(128,11).Borders.LineStyle = 0
$ws.Cells.Item(128,12).Value = "0 hours, -30 minutes"
$ws.Cells.Item(128,13).Borders.LineStyle = 0

# Row 129
$ws.Cells.Item(129,1).Value = 128
$ws.Cells.Item(129,2).Value = "Friday, Jan 13"
$ws.Cells.Item(129,3).Value = "4:45 PM"
$ws.Cells.Item(129,4).Value = "FR9316"
$ws.Cells.Item(129,5).Value = "Helsinki"
$ws.Cells.Item(129,6).Value = "(HEL)"
$ws.Cells.Item(129,7).Value = "Ryanair "
$ws.Cells.Item(129,8).Value = "B738"
$ws.Cells.Item(129,9).Value = "(SP-RKR)"
$ws.Cells.Item(129,10).Value = "4:40 PM"
$ws.Cells.Item(129,11).Borders.LineStyle = 0
$ws.Cells.Item(129,12).Value = "0 hours, -5 minutes"
$ws.Cells.Item(129,13).Borders.LineStyle = 0

# Row 130
$ws.Cells.Item(130,1).Value = 129
$ws.Cells.Item(130,2).Value = "Friday, Jan 13"
$ws.Cells.Item(130,3).Value = "5:00 PM"
$ws.Cells.Item(130,4).Value = "FR7945"
$ws.Cells.Item(130,5).Value = "Leeds"
$ws.Cells.Item(130,6).Value = "(LBA)"
$ws.Cells.Item(130,7).Value = "Ryanair "
$ws.Cells.Item(130,8).Value = "B738"
$ws.Cells.Item(130,9).Value = "(SP-RSX)"
$ws.Cells.Item(130,10).Value = "4:31 PM"
$ws.Cells.Item(130,11).Borders.LineStyle = 0
$ws.Cells.Item(130,12).Value = "0 hours, -29 minutes"
$ws.Cells.Item(130,13).Borders.LineStyle = 0

# Row 131
$ws.Cells.Item(131,1).Value = 130
$ws.Cells.Item(131,2).Value = "Friday, Jan 13"
$ws.Cells.Item(131,3).Value = "5:50 PM"
$ws.Cells.Item(131,4).Value = "LO3947"
$ws.Cells.Item(131,5).Value = "Warsaw"
$ws.Cells.Item(131,6).Value = "(WAW)"
$ws.Cells.Item(131,7).Value = "LOT "
$ws.Cells.Item(131,8).Value = "E170"
$ws.Cells.Item(131,9).Value = "(SP-LDI)"
$ws.Cells.Item(131,10).Value = "5:36 PM"
$ws.Cells.Item(131,11).Borders.LineStyle = 0
$ws.Cells.Item(131,12).Value = "0 hours, -14 minutes"
$ws.Cells.Item(131,13).Borders.LineStyle = 0

# Row 132
$ws.Cells.Item(132,1).Value = 131
$ws.Cells.Item(132,2).Value = "Friday, Jan 13"
$ws.Cells.Item(132,3).Value = "8:55 PM"
$ws.Cells.Item(132,4).Value = "E47014"
$ws.Cells.Item(132,5).Value = "Dubai"
$ws.Cells.Item(132,6).Value = "(DWC)"
$ws.Cells.Item(132,7).Value = "Enter Air "
$ws.Cells.Item(132,8).Value = "B738"
$ws.Cells.Item(132,9).Value = "(SP-ENL)"
$ws.Cells.Item(132,10).Value = "9:41 PM"
$ws.Cells.Item(132,11).Borders.LineStyle = 0
$ws.Cells.Item(132,12).Value = "0 hours, 46 minutes"
$ws.Cells.Item(132,13).Borders.LineStyle = 0

# Row 133
$ws.Cells.Item(133,1).Value = 132
$ws.Cells.Item(133,2).Value = "Friday, Jan 13"
$ws.Cells.Item(133,3).Value = "8:55 PM"
$ws.Cells.Item(133,4).Value = "FR7949"
$ws.Cells.Item(133,5).Value = "Bristol"
$ws.Cells.Item(133,6).Value = "(BRS)"
$ws.Cells.Item(133,7).Value = "Ryanair "
$ws.Cells.Item(133,8).Value = "B738"
$ws.Cells.Item(133,9).Value = "(SP-RSM)"
$ws.Cells.Item(133,10).Value = "8:49 PM"
$ws.Cells.Item(133,11).Borders.LineStyle = 0
$ws.Cells.Item(133,12).Value = "0 hours, -6 minutes"
$ws.Cells.Item(133,13).Borders.LineStyle = 0

# Row 134
$ws.Cells.Item(134,1).Value = 133
$ws.Cells.Item(134,2).Value = "Friday, Jan 13"
$ws.Cells.Item(134,3).Value = "9:20 PM"
$ws.Cells.Item(134,4).Value = "FR7678"
$ws.Cells.Item(134,5).Value = "Stockholm"
$ws.Cells.Item(134,6).Value = "(ARN)"
$ws.Cells.Item(134,7).Value = "Ryanair "
$ws.Cells.Item(134,8).Value = "B38M"
$ws.Cells.Item(134,9).Value = "(9H-VUJ)"
$ws.Cells.Item(134,10).Value = "9:37 PM"
$ws.Cells.Item(134,11).Borders.LineStyle = 0
$ws.Cells.Item(134,12).Value = "0 hours, 17 minutes"
$ws.Cells.Item(134,13).Borders.LineStyle = 0

# Row 135
$ws.Cells.Item(135,1).Value = 134
$ws.Cells.Item(135,2).Value = "Friday, Jan 13"
$ws.Cells.Item(135,3).Value = "9:40 PM"
$ws.Cells.Item(135,4).Value = "FR8320"
$ws.Cells.Item(135,5).Value = "London"
$ws.Cells.Item(135,6).Value = "(STN)"
$ws.Cells.Item(135,7).Value = "Ryanair "
$ws.Cells.Item(135,8).Value = "B38M"
$ws.Cells.Item(135,9).Value = "(EI-HEY)"
$ws.Cells.Item(135,10).Value = "9:47 PM"
$ws.Cells.Item(135,11).Borders.LineStyle = 0
$ws.Cells.Item(135,12).Value = "0 hours, 7 minutes"
$ws.Cells.Item(135,13).Borders.LineStyle = 0

# Row 136
$ws.Cells.Item(136,1).Value = 135
$ws.Cells.Item(136,2).Value = "Friday, Jan 13"
$ws.Cells.Item(136,3).Value = "9:45 PM"
$ws.Cells.Item(136,4).Value = "FR3593"
$ws.Cells.Item(136,5).Value = "Milan"
$ws.Cells.Item(136,6).Value = "(BGY)"
$ws.Cells.Item(136,7).Value = "Malta Air "
$ws.Cells.Item(136,8).Value = "B38M"
$ws.Cells.Item(136,9).Value = "(9H-VUB)"
$ws.Cells.Item(136,10).Value = "9:31 PM"
$ws.Cells.Item(136,11).Borders.LineStyle = 0
$ws.Cells.Item(136,12).Value = "0 hours, -14 minutes"
$ws.Cells.Item(136,13).Borders.LineStyle = 0

# Row 137
$ws.Cells.Item(137,1).Value = 136
$ws.Cells.Item(137,2).Value = "Friday, Jan 13"
$ws.Cells.Item(137,3).Value = "11:00 PM"
$ws.Cells.Item(137,4).Value = "FR1974"
$ws.Cells.Item(137,5).Value = "Dublin"
$ws.Cells.Item(137,6).Value = "(DUB)"
$ws.Cells.Item(137,7).Value = "Ryanair "
$ws.Cells.Item(137,8).Value = "B738"
$ws.Cells.Item(137,9).Value = "(SP-RKR)"
$ws.Cells.Item(137,10).Value = "10:52 PM"
$ws.Cells.Item(137,11).Borders.LineStyle = 0
$ws.Cells.Item(137,12).Value = "0 hours, -8 minutes"
$ws.Cells.Item(137,13).Borders.LineStyle = 0
